$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New right-leaning outlets to add, one per row starting at row 34, column B,
# with a score of 1 in column C (matching existing News Sites (v2) / Score (v2) pattern).
$outlets = @(
    "nationalreview",
    "breitbart",
    "hannity",
    "theblaze",
    "heritage",
    "washingtonexaminer",
    "dailywire",
    "thefederalist",
    "thegatewaypundit",
    "dailycaller",
    "infowars",
    "stanfordreview",
    "thenewamerican",
    "prntly"
)

$startRow = 34
for ($i = 0; $i -lt $outlets.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $outlets[$i]
    $ws.Cells.Item($row, 3).Value = 1
}

# Update the active selection on the sheet to match the new editing location.
$ws.Range("F48").Select()

$wb.Save()
